$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cases" tab was renamed to "Participants" in the automation's tab list.
$ws.Range("A2").Value = "ParticipantsTab"

# Move the saved selection to A2 (previously on B11)
$ws.Range("A2").Select()
